$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 33 (2035_TM152_DBP_Plus_04), shifting the
#        2050-block rows down by one. The inserted row inherits the
#        formatting of row 32 (the row above), which is what we want.
$ws.Rows(33).Insert()

# Row 32 was the previous "current" run for 2035; it no longer is.
$ws.Range("H32").Value = ""

# Fill in the new row 33 with the newest 2035 run.
$ws.Range("A33").Value = "RTP2021"
$ws.Range("B33").Value = 2035
$ws.Range("C33").Value = "2035_TM152_DBP_Plus_04"
$ws.Range("D33").Value = "DraftBlueprint"
$ws.Range("E33").Value = "Plus"
$ws.Range("F33").Value = """Blueprint Plus Crossing (s23)\v1.5.5"""
$ws.Range("G33").Value = "run998"
$ws.Range("H33").Value = "current"

# --- 2. Append a new row after the current last row (44, after the shift
#        above) for the newest 2050 run. Copy formatting from the last
#        existing 2050 row (44) since a brand-new row has no style.
$ws.Range("A44:H44").Copy()
$ws.Range("A45:H45").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A45").Value = "RTP2021"
$ws.Range("B45").Value = 2050
$ws.Range("C45").Value = "2050_TM152_DBP_PlusCrossing_05"
$ws.Range("D45").Value = "DraftBlueprint"
$ws.Range("E45").Value = "Plus"
$ws.Range("F45").Value = """Blueprint Plus Crossing (s23)\v1.5.5"""
$ws.Range("G45").Value = "run998"
$ws.Range("H45").Value = "current"

# --- 3. Column widths: split the old E:F (both width 15) so F becomes
#        wider to fit the new urbansim_path text.
$ws.Range("F1").ColumnWidth = 40.28515625

# --- 4. Cosmetic: update the active selection to reflect where editing
#        ended up.
$ws.Range("H34").Select()
